$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 629
$ws1.Range("F3").Value = 10701
$ws1.Range("F4").Value = 238
$ws1.Range("F8").Value = 12767
$ws1.Range("F9").Value = 13189
$ws1.Range("F11").Value = 1306
$ws1.Range("F12").Value = 5562
$ws1.Range("F13").Value = 929
$ws1.Range("F15").Value = 368
$ws1.Range("F16").Value = 197
$ws1.Range("F17").Value = 1452
$ws1.Range("F18").Value = 371
$ws1.Range("F19").Value = 2043
$ws1.Range("F20").Value = 1065
$ws1.Range("F21").Value = 1615
$ws1.Range("F22").Value = 889
$ws1.Range("F23").Value = 20
$ws1.Range("F24").Value = 522
$ws1.Range("F25").Value = 746
$ws1.Range("F26").Value = 3071
$ws1.Range("F27").Value = 265
$ws1.Range("F28").Value = 2114
$ws1.Range("F29").Value = 14
$ws1.Range("F31").Value = 1706
$ws1.Range("F32").Value = 1018
$ws1.Range("F33").Value = 552
$ws1.Range("F34").Value = 64
$ws1.Range("F35").Value = 112
$ws1.Range("F36").Value = 3825
$ws1.Range("F37").Value = 4489
$ws1.Range("F38").Value = 281
$ws1.Range("F40").Value = 612
$ws1.Range("F42").Value = 3163
$ws1.Range("F45").Value = 308
$ws1.Range("F46").Value = 51
$ws1.Range("F47").Value = 44
$ws1.Range("F48").Value = 4317
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 39
$ws2.Range("F8").Value = 88
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6731
$ws3.Range("F3").Value = 108
$ws3.Range("F4").Value = 178
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 629
$ws4.Range("F3").Value = 10702
$ws4.Range("F5").Value = 108
$ws4.Range("F6").Value = 12768
$ws4.Range("F7").Value = 13189
$ws4.Range("F10").Value = 1306
$ws4.Range("F11").Value = 5562
$ws4.Range("F12").Value = 929
$ws4.Range("F13").Value = 368
$ws4.Range("F14").Value = 88
$ws4.Range("F15").Value = 197
$ws4.Range("F16").Value = 1452
$ws4.Range("F17").Value = 371
$ws4.Range("F18").Value = 2043
$ws4.Range("F19").Value = 1065
$ws4.Range("F20").Value = 1615
$ws4.Range("F21").Value = 889
$ws4.Range("F22").Value = 522
$ws4.Range("F23").Value = 746
$ws4.Range("F24").Value = 3071
$ws4.Range("F26").Value = 265
$ws4.Range("F27").Value = 2114
$ws4.Range("F28").Value = 14
$ws4.Range("F31").Value = 1706
$ws4.Range("F33").Value = 1018
$ws4.Range("F34").Value = 553
$ws4.Range("F35").Value = 64
$ws4.Range("F36").Value = 3825
$ws4.Range("F37").Value = 4489
$ws4.Range("F39").Value = 281
$ws4.Range("F41").Value = 612
$ws4.Range("F43").Value = 3163
$ws4.Range("F45").Value = 308
$ws4.Range("F46").Value = 51
$ws4.Range("F47").Value = 44
$ws4.Range("F48").Value = 4317
